$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.395.49"
$ws.Range("D2").Style = $ws.Range("C2").Style
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "'1.878.38"
$ws.Range("D3").Style = $ws.Range("C3").Style
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'0.7159"
$ws.Range("D5").Style = $ws.Range("C5").Style
$ws.Range("E5").Value = "  +1.10%  "
$ws.Range("D6").Value = "'243.37"
$ws.Range("D6").Style = $ws.Range("C6").Style
$ws.Range("E6").Value = "  +0.60%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "'0.07926"
$ws.Range("D8").Style = $ws.Range("C8").Style
$ws.Range("E8").Value = "  +1.65%  "
$ws.Range("D9").Value = "'0.3141"
$ws.Range("D9").Style = $ws.Range("C9").Style
$ws.Range("E9").Value = "  +1.12%  "
$ws.Range("D10").Value = "'24.92"
$ws.Range("D10").Style = $ws.Range("C10").Style
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("E11").Value = "  -2.89%  "
$ws.Range("D12").Value = "'1.897.80"
$ws.Range("D12").Style = $ws.Range("C12").Style
$ws.Range("E12").Value = "  +1.67%  "
$ws.Range("D13").Value = "'94.79"
$ws.Range("D13").Style = $ws.Range("C13").Style
$ws.Range("E13").Value = "  +3.84%  "
$ws.Range("D14").Value = "'5.232"
$ws.Range("D14").Style = $ws.Range("C14").Style
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("D15").Value = "'0.7066"
$ws.Range("D15").Style = $ws.Range("C15").Style
$ws.Range("E15").Value = "  -1.42%  "
$ws.Range("D16").Value = "'6.382"
$ws.Range("D16").Style = $ws.Range("C16").Style
$ws.Range("E16").Value = "  +3.94%  "
$ws.Range("D17").Value = "'0.000008401"
$ws.Range("D17").Style = $ws.Range("C17").Style
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").Value = "'29.411.92"
$ws.Range("D18").Style = $ws.Range("C18").Style
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("D19").Value = "'252.58"
$ws.Range("D19").Style = $ws.Range("C19").Style
$ws.Range("E19").Value = "  +5.13%  "
$ws.Range("D20").Value = "'13.32"
$ws.Range("D20").Style = $ws.Range("C20").Style
$ws.Range("E20").Value = "  +0.99%  "
$ws.Range("D21").Value = "'2.144.01"
$ws.Range("D21").Style = $ws.Range("C21").Style
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = $ws.Range("C22").Style
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "'7.666"
$ws.Range("D23").Style = $ws.Range("C23").Style
$ws.Range("E23").Value = "  -0.98%  "
$ws.Range("D24").Value = "'1.003"
$ws.Range("D24").Style = $ws.Range("C24").Style
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("D25").Value = "'0.1583"
$ws.Range("D25").Style = $ws.Range("C25").Style
$ws.Range("E25").Value = "  -0.86%  "
$ws.Range("D26").Value = "'9.062"
$ws.Range("D26").Style = $ws.Range("C26").Style
$ws.Range("E26").Value = "  +0.32%  "
$ws.Range("D27").Value = "'162.04"
$ws.Range("D27").Style = $ws.Range("C27").Style
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("D28").Value = "'18.91"
$ws.Range("D28").Style = $ws.Range("C28").Style
$ws.Range("E28").Value = "  +2.41%  "
$ws.Range("D29").Value = "'1.503"
$ws.Range("D29").Style = $ws.Range("C29").Style
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").Value = "'4.408"
$ws.Range("D30").Style = $ws.Range("C30").Style
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").Value = "'4.289"
$ws.Range("D31").Style = $ws.Range("C31").Style
$ws.Range("E31").Value = "  -1.28%  "
$ws.Range("D32").Value = "'1.220"
$ws.Range("D32").Style = $ws.Range("C32").Style
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("D33").Value = "'0.05322"
$ws.Range("D33").Style = $ws.Range("C33").Style
$ws.Range("E33").Value = "  -0.60%  "
$ws.Range("D34").Value = "'1.942"
$ws.Range("D34").Style = $ws.Range("C34").Style
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("D35").Value = "'0.7564"
$ws.Range("D35").Style = $ws.Range("C35").Style
$ws.Range("E35").Value = "  +1.10%  "
$ws.Range("D36").Value = "'1.176"
$ws.Range("D36").Style = $ws.Range("C36").Style
$ws.Range("E36").Value = "  +0.22%  "
$ws.Range("D37").Value = "'2.702"
$ws.Range("D37").Style = $ws.Range("C37").Style
$ws.Range("E37").Value = "  +0.63%  "
$ws.Range("D38").Value = "'0.01891"
$ws.Range("D38").Style = $ws.Range("C38").Style
$ws.Range("E38").Value = "  +0.92%  "
$ws.Range("D39").Value = "'1.274.84"
$ws.Range("D39").Style = $ws.Range("C39").Style
$ws.Range("E39").Value = "  +2.73%  "
$ws.Range("E40").Value = "  +1.22%  "
$ws.Range("D41").Value = "'6.391"
$ws.Range("D41").Style = $ws.Range("C41").Style
$ws.Range("E41").Value = "  -2.04%  "
$ws.Range("D42").Value = "'112.66"
$ws.Range("D42").Style = $ws.Range("C42").Style
$ws.Range("E42").Value = "  +2.58%  "
$ws.Range("D43").Value = "'0.9048"
$ws.Range("D43").Style = $ws.Range("C43").Style
$ws.Range("E43").Value = "  +1.52%  "
$ws.Range("D44").Value = "'73.99"
$ws.Range("D44").Style = $ws.Range("C44").Style
$ws.Range("E44").Value = "  +2.48%  "
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("E46").Value = "  +1.07%  "
$ws.Range("D47").Value = "'2.038.72"
$ws.Range("D47").Style = $ws.Range("C47").Style
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("D48").Value = "'1.805"
$ws.Range("D48").Style = $ws.Range("C48").Style
$ws.Range("E48").Value = "  +0.67%  "
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("D50").Value = "'9.494"
$ws.Range("D50").Style = $ws.Range("C50").Style
$ws.Range("E50").Value = "  +0.60%  "
$ws.Range("D51").Value = "'0.4339"
$ws.Range("D51").Style = $ws.Range("C51").Style
$ws.Range("E51").Value = "  +0.10%  "
